# Commit: "print 1-8 chapter in protocol"
#
# The document currently has chapter "1 Основание для проведения испытаний"
# followed by a separate paragraph with the supporting document reference.
# This edit turns that reference line into a line inside chapter 1 (bracketed
# by <w:br/> line breaks) and appends chapters 2-7 (Информация о заказчике,
# Информация об объекте испытаний, Даты проведения испытаний, Цель испытаний,
# Условия окружающей среды, Методы испытаний) as new paragraphs right after it,
# each built from "<number> <w:br/> label: value <w:br/> ..." runs.

$d = $word.ActiveDocument

# Find the "1 Основание для проведения испытаний" paragraph and the paragraph
# right after it (the "Документы на проведение ..." line) so we can replace
# that pair with the full chapter 1-7 block.
$chapter1 = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text -like "*Основание для проведения испытаний*") {
        $chapter1 = $i
        break
    }
}

$startPara = $d.Paragraphs.Item($chapter1)
$endPara = $d.Paragraphs.Item($chapter1 + 1)
$range = $d.Range($startPara.Range.Start, $endPara.Range.End)

$newXml = '<w:p><w:r><w:rPr><w:b/></w:rPr><w:t>1</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:b/></w:rPr><w:t>Основание для проведения испытаний</w:t></w:r><w:r><w:br/></w:r><w:r><w:t>Документы на проведение сертификационных испытаний ООО «ФаерЛаб» № 0002-2 стр/э/зн от 28.06.2024.</w:t></w:r><w:r><w:br/></w:r></w:p><w:p><w:r><w:rPr><w:b/></w:rPr><w:t>2</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:b/></w:rPr><w:t>Информация о заказчике</w:t></w:r><w:r><w:br/></w:r><w:r><w:t xml:space="preserve">юридический адрес: </w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>143985, Московская область, г. Балашиха, мкр. Железнодорожный, ул. Автозаводская, д. 50а, пом. 16, 16а</w:t></w:r><w:r><w:br/></w:r><w:r><w:t xml:space="preserve">адрес места осуществления деятельности: </w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>143985, Московская область, г. Балашиха, мкр. Железнодорожный, ул. Автозаводская, д. 50в, этаж 1, пом. 33</w:t></w:r><w:r><w:br/></w:r><w:r><w:t xml:space="preserve">наименование: </w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>ООО «ФаерЛаб»</w:t></w:r><w:r><w:br/></w:r><w:r><w:t xml:space="preserve">телефон: </w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>+7 (495) 112-01-93</w:t></w:r><w:r><w:br/></w:r><w:r><w:t xml:space="preserve">e-mail: </w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>info@firelab.su</w:t></w:r><w:r><w:br/></w:r><w:r><w:t xml:space="preserve">номер в реестре аккредитованных лиц: </w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>RA.RU.11НК07</w:t></w:r><w:r><w:br/></w:r></w:p><w:p><w:r><w:rPr><w:b/></w:rPr><w:t>3</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:b/></w:rPr><w:t>Информация об объекте испытаний</w:t></w:r><w:r><w:br/></w:r><w:r><w:t xml:space="preserve">ID: </w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>3926</w:t></w:r><w:r><w:br/></w:r><w:r><w:t xml:space="preserve">Образец представлен на испытания: </w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>19.07.2024</w:t></w:r><w:r><w:br/></w:r><w:r><w:t xml:space="preserve">Марка: </w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>КПпЭфнг(A)-HF 1х2х0,98</w:t></w:r><w:r><w:br/></w:r><w:r><w:t xml:space="preserve">Партия: </w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r/><w:r><w:br/></w:r><w:r><w:t xml:space="preserve">Папка с фото образца: </w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r/><w:r><w:br/></w:r></w:p><w:p><w:r><w:rPr><w:b/></w:rPr><w:t>4</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:b/></w:rPr><w:t>Даты проведения испытаний</w:t></w:r><w:r><w:br/></w:r><w:r><w:t>Дата начала</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>14.07.2024</w:t></w:r><w:r><w:br/></w:r><w:r><w:t>Дата окончания</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>25.07.2024</w:t></w:r><w:r><w:br/></w:r></w:p><w:p><w:r><w:rPr><w:b/></w:rPr><w:t>5</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:b/></w:rPr><w:t>Цель испытаний</w:t></w:r><w:r><w:br/></w:r><w:r><w:t>Определение соответствия кабеля марки КПпЭфнг(A)-HF 1х2х0,98 требованиям ГОСТ 31565-2012"Кабельные изделия. Требования пожарной безопасности" (п. 5.3, 5.5), ТУ 16.К99-027-2005"КАБЕЛИ СИММЕТРИЧНЫЕ ДЛЯ СЕТЕЙ ПРОМЫШЛЕННОЙ АВТОМАТИЗАЦИИ С НИЗКИМ ДЫМО- И ГАЗОВЫДЕЛЕНИЕМ" (п. 1.2.2, 1.2.2, 2.2.1, 1.3.1, 1.3.2, 1.3.3, 1.3.5, 1.3.7, 1.3.8, 1.4.7, 1.5.1, 1.5.2, 1.5.3, 1.6.2, 1.3.14, 1.4.1, табл. 5, п. 1, 1.4.1, табл. 5, п. 2, 1.4.1, табл. 5, п. 3, 1.4.1, табл. 5, п. 4, 1.4.1, табл. 5, п. 5, 1.4.1, табл. 5, п. 6, 1.4.1, табл. 5, п. 7, 1.4.1, табл. 5, п. 7)</w:t></w:r><w:r><w:br/></w:r></w:p><w:p><w:r><w:rPr><w:b/></w:rPr><w:t>6</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:b/></w:rPr><w:t>Условия окружающей среды</w:t></w:r><w:r><w:br/></w:r><w:r><w:t>Температура</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>(21,0 - 24,0) °C</w:t></w:r><w:r><w:br/></w:r><w:r><w:t>Относительная влажность воздуха</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>(44,0 - 56,0) %</w:t></w:r><w:r><w:br/></w:r><w:r><w:t>Атмосферное давление</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>(97,5 - 97,5) кПа</w:t></w:r><w:r><w:br/></w:r></w:p><w:p><w:r><w:rPr><w:b/></w:rPr><w:t>7</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:b/></w:rPr><w:t>Методы испытаний</w:t></w:r><w:r><w:br/></w:r><w:r><w:t>ГОСТ 12177-79</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>Кабели, провода и шнуры. Методы проверки конструкции</w:t></w:r><w:r><w:br/></w:r><w:r><w:t>ГОСТ 20.57.406-81</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>Комплексная система контроля качества. Изделия электронной техники, квантовой электроники и электротехнические методы испытаний</w:t></w:r><w:r><w:br/></w:r><w:r><w:t>ГОСТ 27893-88</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>Кабели связи. Методы испытаний</w:t></w:r><w:r><w:br/></w:r><w:r><w:t>ГОСТ 2990-78</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>Кабели,провода и шнуры. Методы испытания напряжением</w:t></w:r><w:r><w:br/></w:r><w:r><w:t>ГОСТ 3345-76</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>Кабели,провода и шнуры. Метод определения эллектрического сопротивления изоляции</w:t></w:r><w:r><w:br/></w:r><w:r><w:t>ГОСТ 7229-76</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>Кабели,провода и шнуры. Метод определения эллектрического сопротивления токопроводящих жили проводов</w:t></w:r><w:r><w:br/></w:r><w:r><w:t>ГОСТ IEC 60332-3-22-2011</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>Испытания электрических и оптических кабелей в условиях воздействия пламени. Часть 3-22. Распространение пламени по вертикально расположенный пучкам проводов или кабеля. категория A</w:t></w:r><w:r><w:br/></w:r><w:r><w:t>ГОСТ IEC 60811-501-2015</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>Кабели электрические и волоконно-оптические. Методы испытаний неметаллических материалов. Часть 501. Механические испытания испытания для определения механических свойств композиций изоляции и оболочки</w:t></w:r><w:r><w:br/></w:r><w:r><w:t>ГОСТ IEC 61034-2-2011</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>Измерение плотности дыма при горении кабелей в заданных условиях. Часть 2. Метод испытания и требования к нему</w:t></w:r><w:r><w:br/></w:r><w:r><w:t>ГОСТ Р 54429-2011</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>Кабели связи симметричные для цифровых систем передачи. Общие технические условия</w:t></w:r><w:r><w:br/></w:r></w:p>'

$range.InsertXML($newXml)

Write-Output "chapters 1-7 inserted"
